# Generate Report for Handback
# Refresh the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# values for the first data row (the 85b0a023-... file) on each localized
# status sheet, reflecting a newer handback report run.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value = "2016-03-15 03:35:17"
$ws_zhcn.Range("H2").Value = "2016-03-15 03:35:56"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value = "2016-03-15 03:35:25"
$ws_dede.Range("H2").Value = "2016-03-15 03:36:09"
